# Refresh the cryptocurrency snapshot: "Price" (D) and "Volume(1h)" (E) columns,
# driven by the latest GitHub Actions scrape (commit "Updated cryptos list").
#
# Both columns are stored as plain text in the sheet (no numeric cell format),
# so D-column updates are entered with a leading apostrophe. That forces Excel to
# keep them as text -- preserving things like trailing zeros ("94.60", "0.0800")
# and the "."-grouped big numbers ("43.225.76") -- without altering the cell
# NumberFormat (it stays "General", matching the original file) or leaving a
# visible apostrophe in the value. E-column values already read as text (they
# contain "%" and padding spaces), so they are set directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.174.77"
$ws.Range("E2").Value = "  -6.17%  "
$ws.Range("D3").Value = "'2.558.41"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'299.13"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "'94.60"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("E7").Value = "  -3.39%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -5.15%  "
$ws.Range("D10").Value = "'35.99"
$ws.Range("E10").Value = "  -7.65%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D14").Value = "'2.950.51"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "'2.550.02"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "'0.876"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("E17").Value = "  -4.31%  "
$ws.Range("D18").Value = "'43.225.76"
$ws.Range("E18").Value = "  -6.58%  "
$ws.Range("D19").Value = "'13.07"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").Value = "'0.0₃0986"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "'72.33"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "'261.24"
$ws.Range("E23").Value = "  -10.64%  "
$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  -4.10%  "
$ws.Range("D25").Value = "'29.76"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  -7.05%  "
$ws.Range("D29").Value = "'37.14"
$ws.Range("E29").Value = "  -4.66%  "
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("D32").Value = "'154.77"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  -6.60%  "
$ws.Range("D36").Value = "'0.0800"
$ws.Range("E36").Value = "  -5.27%  "
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("D39").Value = "'16.79"
$ws.Range("E39").Value = "  +6.71%  "
$ws.Range("D40").Value = "'23.45"
$ws.Range("E40").Value = "  +8.70%  "
$ws.Range("D41").Value = "'3.50"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("E42").Value = "  -5.22%  "
$ws.Range("D43").Value = "'3.91"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "'2.067.02"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "'86.03"
$ws.Range("E46").Value = "  -10.88%  "
$ws.Range("E47").Value = "  +3.63%  "
$ws.Range("D48").Value = "'2.806.69"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("E49").Value = "  -7.15%  "
$ws.Range("D50").Value = "'1.71"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'104.55"
$ws.Range("E51").Value = "  -4.93%  "
